$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Calculation")

$ws.Range("B7").Value = 554829.3062950152
$ws.Range("C7").Value = 30.6

$ws.Range("B13").Value = 27740.71609267845
$ws.Range("C13").Value = 28.646

$ws.Range("B14").Value = 30816.34669766158
$ws.Range("C14").Value = 31.822

$ws.Range("B15").Value = 87424.63341143535
$ws.Range("C15").Value = 27.646

$ws.Range("B16").Value = 97468.06232392609
$ws.Range("C16").Value = 30.822
